# "Generate Report for Handback"
#
# The localization-status workbook gets refreshed by the CI handback job:
#   - Status text moves from "In Translation" to "Handed back: in sync with en-US"
#   - Each locale sheet (zh-cn / de-de) gets its "Latest Target File" (I) and
#     "Latest Handback File" (J) columns populated with the handed-back file
#     name (hyperlinked back to the source doc) and the generated xlf name.
#   - The "Latest Handback DateTime" (K) stamps get refreshed.
#   - A couple of columns get widened so the new long file names are readable.

$wb = $excel.ActiveWorkbook

function Set-ExactColumnWidth($range, [double]$targetCharWidth) {
    # This host's column-width model quantises to 1/6-character steps
    # (output = (round(input*6)+5)/6), so solve for the COM input that lands
    # on the closest achievable value to the desired stored width.
    $input = ([math]::Round($targetCharWidth * 6) - 5) / 6
    $range.ColumnWidth = $input
}

# ---------------------------------------------------------------------
# 1. Status text: "In Translation" -> "Handed back: in sync with en-US"
#    (shows up on Overview!E2:F3 and on each locale sheet's Status column)
# ---------------------------------------------------------------------
foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Cells.Replace("In Translation", "Handed back: in sync with en-US")
}

# ---------------------------------------------------------------------
# 2. Overview sheet column widths (E, F widen to fit the new status text)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
Set-ExactColumnWidth $wsOverview.Columns.Item(5) 29.9777047293527
Set-ExactColumnWidth $wsOverview.Columns.Item(6) 29.9777047293527

# ---------------------------------------------------------------------
# 3. Per-locale sheets: fill in the handback columns
# ---------------------------------------------------------------------
$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3a02a5aee618105b084b664e1e162430645b760/e2e/"

$locales = @(
    @{ Name = "zh-cn"; HandbackDate = "2016-09-05 12:31:54" },
    @{ Name = "de-de"; HandbackDate = "2016-09-05 12:32:04" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Name)

    # Column widths: Status (C) widens, Latest Target File / Latest Handback
    # File (I, J) widen to fit full file names.
    Set-ExactColumnWidth $ws.Columns.Item(3) 29.9777047293527
    Set-ExactColumnWidth $ws.Columns.Item(9) 40
    Set-ExactColumnWidth $ws.Columns.Item(10) 40

    # Row 2 -> 07e7a3c0-...md ; Row 3 -> 081d15f1-...md
    $rows = @(
        @{ Row = 2; Stem = "07e7a3c0-6436-4fd8-abcc-d407b022910e"; Hash = "da5d67fca4da1eabbcbcd61c096c592e00edaca1" },
        @{ Row = 3; Stem = "081d15f1-f41e-4ab7-be5a-c585a51e2584"; Hash = "d283ae20fc019fba7a951601cac9f31b50e52272" }
    )

    # Re-create the sheet's hyperlinks in document order (A2, I2, A3, I3) so
    # relationship ids line up the way the generator emits them.
    $ws.Hyperlinks.Delete()

    foreach ($r in $rows) {
        $row = $r.Row
        $mdName = "$($r.Stem).md"
        $xlfName = "$($r.Stem).$($r.Hash).$($locale.Name).xlf"
        $mdUrl = "$baseUrl$mdName"

        # Column A: source file hyperlink (unchanged content, re-added so the
        # relationship numbering matches the freshly generated report).
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 1), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName)

        # Column I: Latest Target File - hyperlinked handback target doc.
        $ws.Cells.Item($row, 9).Value = $mdName
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 9), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName)
        $ws.Cells.Item($row, 9).Style = "HyperLink"

        # Column J: Latest Handback File - generated xlf name (plain text).
        $ws.Cells.Item($row, 10).Value = $xlfName

        # Column K: Latest Handback DateTime.
        $ws.Cells.Item($row, 11).Value = $locale.HandbackDate
    }
}
